$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column A ("單位別") entirely, shifting remaining columns left.
$ws.Columns.Item(1).Delete()

Write-Host "Done"
